# Moderhinke Dashboard - weekly data refresh
#
# 1. A handful of already-present weeks get their farms_total_count /
#    farms_to_examine_count figures nudged slightly (a re-pull of the
#    source data before the new week was appended).
# 2. The new ISO week 2025-13 (last day of week 2025-03-30) is appended
#    as five new rows, following the existing farms_total_count /
#    farms_to_examine_count / farms_examined_count /
#    farms_examined_positive_count / farms_examined_negative_count block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Revise existing D-column figures -----------------------------
$revisions = @(
    @{ Cell = "D2";   Value = 11730 },
    @{ Cell = "D3";   Value = 11475 },
    @{ Cell = "D7";   Value = 11834 },
    @{ Cell = "D8";   Value = 11121 },
    @{ Cell = "D12";  Value = 11942 },
    @{ Cell = "D13";  Value = 10619 },
    @{ Cell = "D17";  Value = 11986 },
    @{ Cell = "D18";  Value = 10068 },
    @{ Cell = "D22";  Value = 12044 },
    @{ Cell = "D23";  Value = 9640 },
    @{ Cell = "D27";  Value = 12085 },
    @{ Cell = "D28";  Value = 9180 },
    @{ Cell = "D32";  Value = 12125 },
    @{ Cell = "D33";  Value = 8772 },
    @{ Cell = "D37";  Value = 12163 },
    @{ Cell = "D38";  Value = 8248 },
    @{ Cell = "D42";  Value = 12198 },
    @{ Cell = "D43";  Value = 7681 },
    @{ Cell = "D47";  Value = 12223 },
    @{ Cell = "D48";  Value = 7133 },
    @{ Cell = "D52";  Value = 12251 },
    @{ Cell = "D53";  Value = 6491 },
    @{ Cell = "D57";  Value = 12274 },
    @{ Cell = "D58";  Value = 5880 },
    @{ Cell = "D62";  Value = 12285 },
    @{ Cell = "D63";  Value = 5714 },
    @{ Cell = "D67";  Value = 12307 },
    @{ Cell = "D68";  Value = 5530 },
    @{ Cell = "D72";  Value = 12323 },
    @{ Cell = "D73";  Value = 5094 },
    @{ Cell = "D77";  Value = 12345 },
    @{ Cell = "D78";  Value = 4616 },
    @{ Cell = "D82";  Value = 12363 },
    @{ Cell = "D83";  Value = 4097 },
    @{ Cell = "D87";  Value = 12390 },
    @{ Cell = "D88";  Value = 3691 },
    @{ Cell = "D92";  Value = 12406 },
    @{ Cell = "D93";  Value = 3353 },
    @{ Cell = "D97";  Value = 12427 },
    @{ Cell = "D98";  Value = 3042 },
    @{ Cell = "D102"; Value = 12442 },
    @{ Cell = "D103"; Value = 2766 },
    @{ Cell = "D107"; Value = 12452 },
    @{ Cell = "D108"; Value = 2479 },
    @{ Cell = "D112"; Value = 12473 },
    @{ Cell = "D113"; Value = 2242 },
    @{ Cell = "D117"; Value = 12494 },
    @{ Cell = "D118"; Value = 1997 },
    @{ Cell = "D122"; Value = 12511 },
    @{ Cell = "D123"; Value = 1743 },
    @{ Cell = "D125"; Value = 1467 },
    @{ Cell = "D126"; Value = 9301 }
)

foreach ($rev in $revisions) {
    $ws.Range($rev.Cell).Value = $rev.Value
}

# --- 2. Append the new week (YearWeekIso 202513, LastDayOfWeek 2025-03-30) ---
$newWeekYearWeekIso = 202513
$newWeekLastDay = 45746   # serial date for 2025-03-30
$startRow = 127

$newRows = @(
    @{ Variable = "farms_total_count";            Number = 12546 },
    @{ Variable = "farms_to_examine_count";        Number = 1484 },
    @{ Variable = "farms_examined_count";          Number = 11062 },
    @{ Variable = "farms_examined_positive_count"; Number = 1473 },
    @{ Variable = "farms_examined_negative_count"; Number = 9589 }
)

# Column B carries the same date-number-format (style) as the rest of the
# table; copy that formatting down from the last existing week instead of
# inventing a brand-new style entry.
$ws.Range("B122").Copy()

$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $newWeekYearWeekIso

    $ws.Cells.Item($r, 2).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Cells.Item($r, 2).Value = $newWeekLastDay

    $ws.Cells.Item($r, 3).Value = $row.Variable
    $ws.Cells.Item($r, 4).Value = $row.Number
    $r++
}

$excel.CutCopyMode = $false

# --- 3. Leave the selection where the author last left it ------------
$ws.Range("C102").Select()
